$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns L (param_P_to_demand2) and M (param_Q_to_demand2), mirroring
# the existing K column's header formatting (bold, thin border, centered).
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("K1").Copy($ws.Range("M1"))
$ws.Application.CutCopyMode = $false

$ws.Range("L1").Value = "param_P_to_demand2"
$ws.Range("M1").Value = "param_Q_to_demand2"

# Fill values 500 / 1000 for data rows 2..17
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 12).Value = 500
    $ws.Cells.Item($r, 13).Value = 1000
}
